$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that sits right after the
#    title heading (Heading1) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$metaPara = $titlePara.Next()
if ($metaPara.Range.Text -notmatch "Meta description") {
    throw "Unexpected paragraph after title: $($metaPara.Range.Text)"
}
[void]$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new bold-title paragraph right before the last paragraph
#    (the one that currently holds the italic "Create a Feature Image
#    Prompt..." text).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
if ($lastPara.Range.Text -notmatch "Create a Feature Image Prompt") {
    throw "Unexpected last paragraph: $($lastPara.Range.Text)"
}
$insertAt = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:r/>' +
              '<w:r><w:rPr><w:b/></w:rPr>' +
              '<w:t>Play Book of Sheba Online for Free - Exciting Ancient Egypt Theme</w:t>' +
              '</w:r></w:p>' +
              '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
[void]$insertAt.InsertXML($newParaXml)

# InsertXML needed a trailing sentinel paragraph to force the break; drop it.
$sentinel = $d.Paragraphs.Item($count + 1)
[void]$sentinel.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the italic "Create a Feature Image Prompt..." text with
#    the new meta-description-style sentence (formatting is preserved).
# ------------------------------------------------------------------
$oldText = "Create a Feature Image Prompt: Design an eye-catching, cartoon-style feature image to capture the essence of Book of Sheba. The main focus of the image should be a happy Maya warrior with glasses. The warrior should be dressed in traditional clothing and holding a scepter. The background of the image should feature a desert landscape and ancient Egyptian artifacts like pyramids and hieroglyphs. Use bright, vibrant colors to make the image pop and attract players to this exciting online slot game."
$newText = "Read our review of Book of Sheba, an online slot game with exciting Ancient Egypt theme, numerous paylines, and immersive sound design. Play for free now!"

[void]$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                               $true, 1, $false, $newText, 2)
